# Fruta / hortaliza, semanal
# Weekly update: insert 4 new rows of "Chirimoya" price data (Mercado Mayorista
# Lo Valledor de Santiago, "Cultivar IV Región") for the new reporting date
# 2021-09-10 (serial 44449), shifting the existing rows 50-71 down to 54-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows above the current row 50 (pushes old rows 50-71 -> 54-75)
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()

$newRows = @(
  @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",44449,13,"Fruta",100107,"Otros",100107002,"Chirimoya","Cultivar IV Región","Especial",125,3200,3200,3200,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3200,1),
  @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",44449,13,"Fruta",100107,"Otros",100107002,"Chirimoya","Cultivar IV Región","Extra (doble especial)",85,3500,3500,3500,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",3500,1),
  @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",44449,13,"Fruta",100107,"Otros",100107002,"Chirimoya","Cultivar IV Región","Primera",150,2700,3000,2850,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",2850,1),
  @(6,"Mercado Mayorista Lo Valledor de Santiago","Metropolitana",44449,13,"Fruta",100107,"Otros",100107002,"Chirimoya","Cultivar IV Región","Segunda",125,2200,2200,2200,"`$/kilo (en caja de 15 kilos)","Provincia de Limarí",2200,1)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $targetRow = 50 + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($targetRow, 1 + $j).Value = $rowData[$j]
    }
}
